# Recreated mantel correlograms with Euclidean distances.
# Update the "Mantel r" and "p" columns of the single results table
# (header row = 1, so data rows are 2..15; col 3 = Mantel r, col 4 = p).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{ Row = 2;  Col = 3; Old = "-0.020"; New = "0.008" },
    @{ Row = 2;  Col = 4; Old = "0.07";   New = "0.257" },
    @{ Row = 3;  Col = 4; Old = "0.485";  New = "0.513" },
    @{ Row = 4;  Col = 3; Old = "-0.016"; New = "0.012" },
    @{ Row = 4;  Col = 4; Old = "0.428";  New = "0.77" },
    @{ Row = 5;  Col = 3; Old = "-0.002"; New = "0.008" },
    @{ Row = 5;  Col = 4; Old = "0.937";  New = "1" },
    @{ Row = 6;  Col = 3; Old = "0.003";  New = "-0.004" },
    @{ Row = 7;  Col = 3; Old = "0.004";  New = "-0.005" },
    @{ Row = 8;  Col = 3; Old = "-0.005"; New = "-0.008" },
    @{ Row = 9;  Col = 3; Old = "-0.007"; New = "0.002" },
    @{ Row = 10; Col = 3; Old = "0.006";  New = "0.032" },
    @{ Row = 10; Col = 4; Old = "1";      New = "0.279" },
    @{ Row = 11; Col = 3; Old = "-0.006"; New = "0.003" },
    @{ Row = 12; Col = 3; Old = "0.009";  New = "0.014" },
    @{ Row = 13; Col = 3; Old = "0.029";  New = "-0.014" },
    @{ Row = 13; Col = 4; Old = "0.503";  New = "1" },
    @{ Row = 14; Col = 3; Old = "0.001";  New = "-0.013" },
    @{ Row = 15; Col = 3; Old = "-0.002"; New = "-0.001" }
)

foreach ($c in $changes) {
    $cell = $t.Cell($c.Row, $c.Col)
    $rng = $cell.Range
    $current = $rng.Text.TrimEnd([char]7, [char]13)
    if ($current -ne $c.Old) {
        Write-Host "MISMATCH at row" $c.Row "col" $c.Col "expected" $c.Old "got" $current
    }
    $rng.Text = $c.New
}
